{"js": "// Replace the date line and all \"A\u00f7B=\" exercise values with their new\n// values, as described by the diff. Each old value is a unique literal\n// string in the document, so a simple search + insertText(\"...\", \"Replace\")\n// per pair is sufficient and keeps existing run formatting intact.\nconst replacements = [\n  [\"2025-07-22 Tuesday\", \"2025-07-23 Wednesday\"],\n  [\"114\u00f77=\", \"853\u00f77=\"],\n  [\"967\u00f75=\", \"778\u00f78=\"],\n  [\"331\u00f75=\", \"187\u00f74=\"],\n  [\"554\u00f79=\", \"282\u00f76=\"],\n  [\"915\u00f76=\", \"959\u00f78=\"],\n  [\"362\u00f75=\", \"822\u00f74=\"],\n  [\"808\u00f78=\", \"488\u00f72=\"],\n  [\"955\u00f78=\", \"603\u00f73=\"],\n  [\"576\u00f75=\", \"720\u00f74=\"],\n  [\"920\u00f77=\", \"627\u00f73=\"],\n  [\"201\u00f78=\", \"689\u00f74=\"],\n  [\"420\u00f78=\", \"445\u00f78=\"],\n  [\"455\u00f73=\", \"779\u00f75=\"],\n  [\"123\u00f78=\", \"274\u00f72=\"],\n  [\"650\u00f72=\", \"798\u00f73=\"],\n  [\"330\u00f72=\", \"300\u00f75=\"],\n  [\"182\u00f72=\", \"437\u00f76=\"],\n  [\"782\u00f76=\", \"231\u00f78=\"],\n  [\"695\u00f77=\", \"661\u00f76=\"],\n  [\"125\u00f74=\", \"687\u00f76=\"],\n  [\"739\u00f78=\", \"922\u00f77=\"],\n  [\"948\u00f73=\", \"608\u00f72=\"],\n  [\"571\u00f75=\", \"343\u00f73=\"],\n  [\"377\u00f79=\", \"957\u00f75=\"],\n  [\"315\u00f78=\", \"954\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and all \"A\u00f7B=\" exercise values with their new\n# values, as described by the diff. Each old value is a unique literal\n# string in the document, so Find/Replace (one pair at a time, over the\n# whole document content range) is sufficient and preserves existing run\n# formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-07-22 Tuesday\", \"2025-07-23 Wednesday\"),\n    @(\"114\u00f77=\", \"853\u00f77=\"),\n    @(\"967\u00f75=\", \"778\u00f78=\"),\n    @(\"331\u00f75=\", \"187\u00f74=\"),\n    @(\"554\u00f79=\", \"282\u00f76=\"),\n    @(\"915\u00f76=\", \"959\u00f78=\"),\n    @(\"362\u00f75=\", \"822\u00f74=\"),\n    @(\"808\u00f78=\", \"488\u00f72=\"),\n    @(\"955\u00f78=\", \"603\u00f73=\"),\n    @(\"576\u00f75=\", \"720\u00f74=\"),\n    @(\"920\u00f77=\", \"627\u00f73=\"),\n    @(\"201\u00f78=\", \"689\u00f74=\"),\n    @(\"420\u00f78=\", \"445\u00f78=\"),\n    @(\"455\u00f73=\", \"779\u00f75=\"),\n    @(\"123\u00f78=\", \"274\u00f72=\"),\n    @(\"650\u00f72=\", \"798\u00f73=\"),\n    @(\"330\u00f72=\", \"300\u00f75=\"),\n    @(\"182\u00f72=\", \"437\u00f76=\"),\n    @(\"782\u00f76=\", \"231\u00f78=\"),\n    @(\"695\u00f77=\", \"661\u00f76=\"),\n    @(\"125\u00f74=\", \"687\u00f76=\"),\n    @(\"739\u00f78=\", \"922\u00f77=\"),\n    @(\"948\u00f73=\", \"608\u00f72=\"),\n    @(\"571\u00f75=\", \"343\u00f73=\"),\n    @(\"377\u00f79=\", \"957\u00f75=\"),\n    @(\"315\u00f78=\", \"954\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
